# Auto-generated Excel COM-interop script to apply scheduled runner updates
# to the Hyperion_Profits workbook (market price / leve profit recalculation).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 17
$ws.Range("H17").Value = 5115.5
$ws.Range("J17").Value = 5115.5
$ws.Range("L17").Value = 15346.5
$ws.Range("N17").Value = -15682.5
# Row 48
$ws.Range("H48").Value = 949.5
$ws.Range("J48").Value = 949.5
$ws.Range("L48").Value = 2848.5
$ws.Range("N48").Value = -3432.5
# Row 56
$ws.Range("H56").Value = 949.5
$ws.Range("J56").Value = 949.5
$ws.Range("L56").Value = 2848.5
$ws.Range("N56").Value = -3916.5
# Row 64
$ws.Range("H64").Value = 7916.8857
$ws.Range("I64").Value = 7249.75
$ws.Range("J64").Value = 8114.5557
$ws.Range("K64").Value = 7249.75
$ws.Range("L64").Value = 8114.5557
$ws.Range("M64").Value = -7001.75
$ws.Range("N64").Value = -8610.555700000001
# Row 67
$ws.Range("H67").Value = 7916.8857
$ws.Range("I67").Value = 7249.75
$ws.Range("J67").Value = 8114.5557
$ws.Range("K67").Value = 7249.75
$ws.Range("L67").Value = 8114.5557
$ws.Range("M67").Value = -6391.75
$ws.Range("N67").Value = -9830.555700000001
# Row 100
$ws.Range("H100").Value = 3797.5557
$ws.Range("I100").Value = 3797.5557
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 3797.5557
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = -3256.5557
$ws.Range("N100").ClearContents()
# Row 132
$ws.Range("H132").Value = 21740788
$ws.Range("I132").Value = 25642738
$ws.Range("K132").Value = 76928214
$ws.Range("M132").Value = -76925684
# Row 138
$ws.Range("H138").Value = 3144.2886
$ws.Range("I138").Value = 1969.125
$ws.Range("J138").Value = 3722.8308
$ws.Range("K138").Value = 5907.375
$ws.Range("L138").Value = 11168.4924
$ws.Range("M138").Value = -767.375
$ws.Range("N138").Value = -21448.4924

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 1984730.9
$ws.Range("I2").Value = 1984730.9
$ws.Range("K2").Value = 1984730.9
$ws.Range("M2").Value = -1984617.9
# Row 32
$ws.Range("H32").Value = 5044.4243
$ws.Range("I32").Value = 3163.64
$ws.Range("J32").Value = 10921.875
$ws.Range("K32").Value = 3163.64
$ws.Range("L32").Value = 10921.875
$ws.Range("M32").Value = -2876.64
$ws.Range("N32").Value = -11495.875
# Row 61
$ws.Range("H61").Value = 4239.7803
$ws.Range("I61").Value = 4445.722
$ws.Range("J61").Value = 2757
$ws.Range("K61").Value = 4445.722
$ws.Range("L61").Value = 2757
$ws.Range("M61").Value = -4233.722
$ws.Range("N61").Value = -3181
# Row 74
$ws.Range("H74").Value = 127315.664
$ws.Range("I74").Value = 84571.586
$ws.Range("J74").Value = 298292
$ws.Range("K74").Value = 84571.586
$ws.Range("L74").Value = 298292
$ws.Range("M74").Value = -83697.586
$ws.Range("N74").Value = -300040
# Row 77
$ws.Range("H77").Value = 127315.664
$ws.Range("I77").Value = 84571.586
$ws.Range("J77").Value = 298292
$ws.Range("K77").Value = 422857.93
$ws.Range("L77").Value = 1491460
$ws.Range("M77").Value = -418489.93
$ws.Range("N77").Value = -1500196
# Row 88
$ws.Range("H88").Value = 1700
$ws.Range("I88").Value = 1600
$ws.Range("J88").Value = 1800
$ws.Range("K88").Value = 1600
$ws.Range("L88").Value = 1800
$ws.Range("M88").Value = -1194
$ws.Range("N88").Value = -2612
# Row 91
$ws.Range("H91").Value = 1700
$ws.Range("I91").Value = 1600
$ws.Range("J91").Value = 1800
$ws.Range("K91").Value = 1600
$ws.Range("L91").Value = 1800
$ws.Range("M91").Value = -196
$ws.Range("N91").Value = -4608
# Row 102
$ws.Range("H102").Value = 20840262
$ws.Range("I102").Value = 41671668
$ws.Range("J102").Value = 8855.5
$ws.Range("K102").Value = 41671668
$ws.Range("L102").Value = 8855.5
$ws.Range("M102").Value = -41670046
$ws.Range("N102").Value = -12099.5
# Row 116
$ws.Range("H116").Value = 1984730.9
$ws.Range("I116").Value = 1984730.9
$ws.Range("K116").Value = 1984730.9
$ws.Range("M116").Value = -1982436.9
# Row 136
$ws.Range("H136").Value = 4239.7803
$ws.Range("I136").Value = 4445.722
$ws.Range("J136").Value = 2757
$ws.Range("K136").Value = 13337.166
$ws.Range("L136").Value = 8271
$ws.Range("M136").Value = -10787.166
$ws.Range("N136").Value = -13371

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 1984730.9
$ws.Range("I3").Value = 1984730.9
$ws.Range("K3").Value = 1984730.9
$ws.Range("M3").Value = -1984616.9
# Row 20
$ws.Range("H20").Value = 2346.1333
$ws.Range("I20").Value = 1899
$ws.Range("J20").Value = 4134.6665
$ws.Range("K20").Value = 1899
$ws.Range("L20").Value = 4134.6665
$ws.Range("M20").Value = -1652
$ws.Range("N20").Value = -4628.6665
# Row 86
$ws.Range("H86").Value = 4555885.5
$ws.Range("I86").Value = 4769499
$ws.Range("K86").Value = 4769499
$ws.Range("M86").Value = -4768376
# Row 89
$ws.Range("H89").Value = 4555885.5
$ws.Range("I89").Value = 4769499
$ws.Range("K89").Value = 23847495
$ws.Range("M89").Value = -23841879
# Row 94
$ws.Range("H94").Value = 3578835.2
$ws.Range("I94").Value = 4000965.5
$ws.Range("J94").Value = 61083.332
$ws.Range("K94").Value = 4000965.5
$ws.Range("L94").Value = 61083.332
$ws.Range("M94").Value = -4000514.5
$ws.Range("N94").Value = -61985.332
# Row 99
$ws.Range("H99").Value = 13080784
$ws.Range("I99").Value = 15986459
$ws.Range("K99").Value = 15986459
$ws.Range("M99").Value = -15984961

$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 1920.2858
$ws.Range("I16").Value = 1362.6364
$ws.Range("J16").Value = 3965
$ws.Range("K16").Value = 1362.6364
$ws.Range("L16").Value = 3965
$ws.Range("M16").Value = -1075.6364
$ws.Range("N16").Value = -4539
# Row 31
$ws.Range("H31").Value = 5015.9707
$ws.Range("I31").Value = 9948.846
$ws.Range("J31").Value = 1962.2858
$ws.Range("K31").Value = 9948.846
$ws.Range("L31").Value = 1962.2858
$ws.Range("M31").Value = -9653.846
$ws.Range("N31").Value = -2552.2858
# Row 32
$ws.Range("H32").Value = 2036.6666
$ws.Range("I32").Value = 2036.6666
$ws.Range("K32").Value = 2036.6666
$ws.Range("M32").Value = -1720.6666
# Row 34
$ws.Range("H34").Value = 5015.9707
$ws.Range("I34").Value = 9948.846
$ws.Range("J34").Value = 1962.2858
$ws.Range("K34").Value = 9948.846
$ws.Range("L34").Value = 1962.2858
$ws.Range("M34").Value = -9746.846
$ws.Range("N34").Value = -2366.2858
# Row 106
$ws.Range("H106").Value = 29995.5
$ws.Range("J106").Value = 29995.5
$ws.Range("L106").Value = 29995.5
$ws.Range("N106").Value = -32519.5
# Row 113
$ws.Range("H113").Value = 1920.2858
$ws.Range("I113").Value = 1362.6364
$ws.Range("J113").Value = 3965
$ws.Range("K113").Value = 1362.6364
$ws.Range("L113").Value = 3965
$ws.Range("M113").Value = 807.3635999999999
$ws.Range("N113").Value = -8305
# Row 132
$ws.Range("H132").Value = 20935.596
$ws.Range("J132").Value = 2825.6
$ws.Range("L132").Value = 8476.799999999999
$ws.Range("N132").Value = -13536.8
# Row 134
$ws.Range("H134").Value = 11263.733
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 11263.733
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 33791.199
$ws.Range("M134").ClearContents()
$ws.Range("N134").Value = -38861.199

$ws = $wb.Worksheets.Item("CUL")
# Row 12
$ws.Range("H12").Value = 109.07692
$ws.Range("J12").Value = 123.22222
$ws.Range("L12").Value = 369.66666
$ws.Range("N12").Value = -715.66666
# Row 138
$ws.Range("H138").Value = 3488.125
$ws.Range("I138").Value = 3488.125
$ws.Range("K138").Value = 10464.375
$ws.Range("M138").Value = -5324.375

$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 15390540
$ws.Range("I70").Value = 22228138
$ws.Range("K70").Value = 22228138
$ws.Range("M70").Value = -22227868
# Row 73
$ws.Range("H73").Value = 15390540
$ws.Range("I73").Value = 22228138
$ws.Range("K73").Value = 22228138
$ws.Range("M73").Value = -22227202
# Row 80
$ws.Range("H80").Value = 1635051.1
$ws.Range("I80").Value = 3500242
$ws.Range("K80").Value = 3500242
$ws.Range("M80").Value = -3499244
# Row 83
$ws.Range("H83").Value = 1635051.1
$ws.Range("I83").Value = 3500242
$ws.Range("K83").Value = 17501210
$ws.Range("M83").Value = -17496218

$ws = $wb.Worksheets.Item("LTW")
# Row 82
$ws.Range("H82").Value = 1985102
$ws.Range("I82").Value = 5051264.5
$ws.Range("K82").Value = 5051264.5
$ws.Range("M82").Value = -5050903.5
# Row 85
$ws.Range("H85").Value = 1985102
$ws.Range("I85").Value = 5051264.5
$ws.Range("K85").Value = 5051264.5
$ws.Range("M85").Value = -5050016.5
# Row 132
$ws.Range("H132").Value = 8734.094999999999
$ws.Range("I132").Value = 8865.511
$ws.Range("J132").Value = 7994.875
$ws.Range("K132").Value = 26596.533
$ws.Range("L132").Value = 23984.625
$ws.Range("M132").Value = -24066.533
$ws.Range("N132").Value = -29044.625
# Row 136
$ws.Range("H136").Value = 36959.414
$ws.Range("I136").Value = 46079
$ws.Range("K136").Value = 138237
$ws.Range("M136").Value = -135687

$ws = $wb.Worksheets.Item("WVR")
# Row 81
$ws.Range("H81").Value = 11906326
$ws.Range("I81").Value = 12821936
$ws.Range("J81").Value = 3400
$ws.Range("K81").Value = 25643872
$ws.Range("L81").Value = 6800
$ws.Range("M81").Value = -25642811
$ws.Range("N81").Value = -8922
# Row 84
$ws.Range("H84").Value = 11906326
$ws.Range("I84").Value = 12821936
$ws.Range("J84").Value = 3400
$ws.Range("K84").Value = 128219360
$ws.Range("L84").Value = 34000
$ws.Range("M84").Value = -128214056
$ws.Range("N84").Value = -44608
# Row 107
$ws.Range("H107").Value = 35719936
$ws.Range("I107").Value = 45457444
$ws.Range("K107").Value = 136372332
$ws.Range("M107").Value = -136370412
